# fix(backend): fix export column's order
#
# The "Colis" (parcel) block and the "Courriers" (mail) block on the
# "Stats structure" sheet were in the wrong order: rows 139-141 held the
# "Colis ..." labels and rows 142-144 held the "Courriers ..." labels.
# Swap the two 3-row blocks so "Courriers ..." comes first (139-141) and
# "Colis ..." comes second (142-144) - the label text stays associated
# with its own row-position semantics but the rows swap content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stats structure")

$ws.Range("B139").Value = "Courriers enregistrés"
$ws.Range("B140").Value = "Courriers remis"
$ws.Range("B141").Value = "Courriers réexpédiés"
$ws.Range("B142").Value = "Colis enregistrés"
$ws.Range("B143").Value = "Colis remis"
$ws.Range("B144").Value = "Colis rééxpédiés"

$ws.Range("B142").Select()
